$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.615.74"
$ws.Range("E2").Value = "  +0.83%  "

$ws.Range("D3").Value = "'1.641.08"
$ws.Range("E3").Value = "  +1.13%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'214.61"
$ws.Range("E5").Value = "  +1.01%  "

$ws.Range("D6").Value = "'0.504"
$ws.Range("E6").Value = "  +1.54%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  +1.14%  "

$ws.Range("D9").Value = "'0.0625"
$ws.Range("E9").Value = "  +0.65%  "

$ws.Range("D10").Value = "'19.10"
$ws.Range("E10").Value = "  +0.96%  "

$ws.Range("D11").Value = "'0.0842"
$ws.Range("E11").Value = "  +0.16%  "

$ws.Range("D12").Value = "'1.868.50"
$ws.Range("E12").Value = "  +1.07%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.18"
$ws.Range("E13").Value = "  +1.83%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.623.74"
$ws.Range("E14").Value = "  -0.34%  "

$ws.Range("D15").Value = "'0.529"
$ws.Range("E15").Value = "  +1.60%  "

$ws.Range("D16").Value = "'64.84"
$ws.Range("E16").Value = "  +1.19%  "

$ws.Range("D17").Value = "'26.611.30"
$ws.Range("E17").Value = "  +0.79%  "

$ws.Range("E18").Value = "  +0.36%  "

$ws.Range("D19").Value = "'215.47"
$ws.Range("E19").Value = "  +0.64%  "

$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").Value = "'4.34"
$ws.Range("E21").Value = "  +1.00%  "

$ws.Range("D22").Value = "'6.25"
$ws.Range("E22").Value = "  +0.58%  "

$ws.Range("D23").Value = "'9.45"
$ws.Range("E23").Value = "  +1.83%  "

$ws.Range("D24").Value = "'2.20"
$ws.Range("E24").Value = "  +12.29%  "

$ws.Range("D25").Value = "'145.01"
$ws.Range("E25").Value = "  -1.69%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("E28").Value = "  +4.32%  "

$ws.Range("D29").Value = "'15.68"
$ws.Range("E29").Value = "  +0.76%  "

$ws.Range("E30").Value = "  +0.86%  "

$ws.Range("E31").Value = "  +1.70%  "

$ws.Range("E32").Value = "  +1.39%  "

$ws.Range("D33").Value = "'3.00"
$ws.Range("E33").Value = "  +1.76%  "

$ws.Range("D34").Value = "'1.277.46"
$ws.Range("E34").Value = "  +5.41%  "

$ws.Range("D35").Value = "'1.54"
$ws.Range("E35").Value = "  +3.17%  "

$ws.Range("E36").Value = "  +1.32%  "

$ws.Range("E37").Value = "  +2.89%  "

$ws.Range("E38").Value = "  +6.09%  "

$ws.Range("D39").Value = "'0.825"
$ws.Range("E39").Value = "  +3.52%  "

$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("D41").Value = "'0.810"
$ws.Range("E41").Value = "  +2.36%  "

$ws.Range("D42").Value = "'2.25"
$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("D43").Value = "'5.41"
$ws.Range("E43").Value = "  +0.81%  "

$ws.Range("D44").Value = "'1.779.00"
$ws.Range("E44").Value = "  +1.10%  "

$ws.Range("D45").Value = "'91.41"
$ws.Range("E45").Value = "  -1.28%  "

$ws.Range("D46").Value = "'59.16"
$ws.Range("E46").Value = "  +8.31%  "

$ws.Range("E47").Value = "  +1.04%  "

$ws.Range("D48").Value = "'0.0515"
$ws.Range("E48").Value = "  +0.96%  "

$ws.Range("D49").Value = "'7.71"
$ws.Range("E49").Value = "  +1.64%  "

$ws.Range("D50").Value = "'0.0962"
$ws.Range("E50").Value = "  +1.32%  "

$ws.Range("D51").Value = "'0.406"
$ws.Range("E51").Value = "  -0.44%  "
